$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of OHLC data (rows 115-117) as produced by the R script.
# Columns: date, volume, high, low, open, close, adj_close(text), ticker
$newRows = @(
    @(45454.2916666667, 0,     1.52999997138977, 1.52999997138977, 1.52999997138977, 1.52999997138977, "1.52999997138977", "SMN.MI"),
    @(45455.2916666667, 0,     1.52999997138977, 1.52999997138977, 1.52999997138977, 1.52999997138977, "1.52999997138977", "SMN.MI"),
    @(45456.5304513889, 15000, 1.54499995708466, 1.5,              1.52999997138977, 1.54499995708466, "1.54499995708466", "SMN.MI")
)

$templateRow = 114
$startRow = 115

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Reuse the date cell's existing style (yyyy-mm-dd hh:mm:ss) instead of
    # creating a brand-new style entry.
    $ws.Range("A$templateRow").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # adj_close is stored as text (it mirrors the shared-string numeric text
    # used elsewhere in this column), so force text formatting, then drop
    # back to the default style so no extra style entry is introduced.
    $adjCell = $ws.Cells.Item($r, 7)
    $adjCell.NumberFormat = "@"
    $adjCell.Value = $row[6]
    $adjCell.Style = "Normal"

    $ws.Cells.Item($r, 8).Value = $row[7]
}
